# Apply cryptos list update (price/volume refresh + one ranking swap)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D holding numeric-looking text (e.g. "247.15") must be
# marked as Text before assignment, otherwise Excel auto-converts the
# string into a floating point number and the "26.524.07"-style grouped-
# thousand text formatting used throughout this sheet would be lost.
$ws.Range("D5:D9").NumberFormat = "@"
$ws.Range("D11:D17").NumberFormat = "@"
$ws.Range("D19:D20").NumberFormat = "@"
$ws.Range("D22:D26").NumberFormat = "@"
$ws.Range("D28:D35").NumberFormat = "@"
$ws.Range("D37:D40").NumberFormat = "@"
$ws.Range("D42:D43").NumberFormat = "@"
$ws.Range("D45:D46").NumberFormat = "@"
$ws.Range("D48:D51").NumberFormat = "@"

# Row-by-row value updates (sheet row number == cell row)
$ws.Range("D2").Value = "26.545.70"
$ws.Range("E2").Value = "  +0.20%  "
$ws.Range("D3").Value = "1.739.79"
$ws.Range("E3").Value = "  +0.29%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "247.15"
$ws.Range("E5").Value = "  +1.25%  "
$ws.Range("D6").Value = "0.9997"
$ws.Range("E6").Value = "  -0.06%  "
$ws.Range("D7").Value = "0.4906"
$ws.Range("E7").Value = "  +2.40%  "
$ws.Range("D8").Value = "0.2670"
$ws.Range("E8").Value = "  +0.22%  "
$ws.Range("D9").Value = "0.06305"
$ws.Range("E9").Value = "  +1.30%  "
$ws.Range("D10").Value = "1.733.66"
$ws.Range("E10").Value = "  -0.08%  "
$ws.Range("D11").Value = "0.07046"
$ws.Range("E11").Value = "  -1.09%  "
$ws.Range("D12").Value = "15.73"
$ws.Range("E12").Value = "  +0.41%  "
$ws.Range("D13").Value = "4.614"
$ws.Range("E13").Value = "  +2.05%  "
$ws.Range("D14").Value = "0.6118"
$ws.Range("E14").Value = "  -0.19%  "
$ws.Range("D15").Value = "77.49"
$ws.Range("E15").Value = "  +0.80%  "
$ws.Range("D16").Value = "0.9995"
$ws.Range("E16").Value = "  -0.07%  "
$ws.Range("D17").Value = "0.000007428"
$ws.Range("E17").Value = "  +7.73%  "
$ws.Range("D18").Value = "26.530.58"
$ws.Range("E18").Value = "  +0.10%  "
$ws.Range("D19").Value = "0.9994"
$ws.Range("E19").Value = "  -0.10%  "
$ws.Range("D20").Value = "11.55"
$ws.Range("E20").Value = "  -1.28%  "
$ws.Range("D21").Value = "1.957.06"
$ws.Range("E21").Value = "  -0.10%  "
$ws.Range("D22").Value = "4.581"
$ws.Range("E22").Value = "  +0.39%  "
$ws.Range("D23").Value = "8.721"
$ws.Range("E23").Value = "  -1.82%  "
$ws.Range("D24").Value = "5.256"
$ws.Range("E24").Value = "  -1.36%  "
$ws.Range("D25").Value = "140.63"
$ws.Range("E25").Value = "  +3.39%  "
$ws.Range("D26").Value = "15.48"
$ws.Range("E26").Value = "  +0.87%  "
$ws.Range("E27").Value = "  +1.25%  "
$ws.Range("D28").Value = "1.770"
$ws.Range("E28").Value = "  -1.33%  "
$ws.Range("D29").Value = "108.04"
$ws.Range("E29").Value = "  +1.33%  "
$ws.Range("D30").Value = "4.052"
$ws.Range("E30").Value = "  +1.80%  "
$ws.Range("D31").Value = "0.08051"
$ws.Range("E31").Value = "  +0.96%  "
$ws.Range("D32").Value = "3.726"
$ws.Range("E32").Value = "  +0.33%  "
$ws.Range("D33").Value = "0.04592"
$ws.Range("E33").Value = "  +0.98%  "
$ws.Range("D34").Value = "0.9992"
$ws.Range("E34").Value = "  -0.06%  "
$ws.Range("D35").Value = "2.609"
$ws.Range("E36").Value = "  +1.95%  "
$ws.Range("D37").Value = "0.6379"
$ws.Range("E37").Value = "  +0.32%  "
$ws.Range("D38").Value = "0.8941"
$ws.Range("E38").Value = "  -4.07%  "
$ws.Range("D39").Value = "2.026"
$ws.Range("E39").Value = "  +2.30%  "
$ws.Range("D40").Value = "2.402"
$ws.Range("E40").Value = "  -0.25%  "
$ws.Range("E41").Value = "  -0.02%  "
$ws.Range("D42").Value = "0.01507"
$ws.Range("E42").Value = "  +0.08%  "
$ws.Range("D43").Value = "102.53"
$ws.Range("E43").Value = "  -6.79%  "
$ws.Range("E44").Value = "  -5.29%  "
$ws.Range("D45").Value = "0.3906"
$ws.Range("E45").Value = "  +0.20%  "
$ws.Range("D46").Value = "6.905"
$ws.Range("E46").Value = "  +0.16%  "
$ws.Range("E47").Value = "  -0.20%  "
$ws.Range("D48").Value = "0.05400"
$ws.Range("E48").Value = "  +1.17%  "
$ws.Range("B49").Value = "Elrond"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D49").Value = "30.61"
$ws.Range("E49").Value = "  -0.49%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "7.790"
$ws.Range("E50").Value = "  -1.47%  "
$ws.Range("D51").Value = "1.277"
$ws.Range("E51").Value = "  +1.03%  "
